$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 68.755
$ws.Range("D2").Value = 68.755
$ws.Range("E2").Value = 2.47443577
$ws.Range("F2").Value = 0.02845953
$ws.Range("G2").Value = 1.95504452
$ws.Range("H2").Value = 135.14081424
$ws.Range("I2").Value = 5.538198863046238
$ws.Range("J2").Value = 5.538198863046238
$ws.Range("K2").Value = 0.2055555307887505
$ws.Range("L2").Value = 0.005253220727049941
$ws.Range("M2").Value = 0.3843323614845019
$ws.Range("N2").Value = 31.5375894163897

$ws.Range("C3").Value = 88.092
$ws.Range("D3").Value = 88.092
$ws.Range("E3").Value = 1.95103218
$ws.Range("F3").Value = 0.02192588
$ws.Range("G3").Value = 1.91246521
$ws.Range("H3").Value = 169.68301096
$ws.Range("I3").Value = 11.59042164348967
$ws.Range("J3").Value = 11.59042164348967
$ws.Range("K3").Value = 0.2519866775246213
$ws.Range("L3").Value = 0.004367517857758321
$ws.Range("M3").Value = 0.3702913988772064
$ws.Range("N3").Value = 45.64084874438079

$ws.Range("C4").Value = 37.692
$ws.Range("D4").Value = 75.317
$ws.Range("E4").Value = 2.27866424
$ws.Range("F4").Value = 0.05143636
$ws.Range("G4").Value = 0.9682347200000001
$ws.Range("H4").Value = 37.03742679
$ws.Range("I4").Value = 4.840099437259792
$ws.Range("J4").Value = 9.668505106322847
$ws.Range("K4").Value = 0.2755381434231993
$ws.Range("L4").Value = 0.01041610135885389
$ws.Range("M4").Value = 0.2275413287453706
$ws.Range("N4").Value = 12.10297667457231

$ws.Range("C5").Value = 47.74
$ws.Range("D5").Value = 93.778
$ws.Range("E5").Value = 1.84017406
$ws.Range("F5").Value = 0.03628156999999999
$ws.Range("G5").Value = 0.85726444
$ws.Range("H5").Value = 41.48008469000001
$ws.Range("I5").Value = 7.436242175837047
$ws.Range("J5").Value = 13.7162988029881
$ws.Range("K5").Value = 0.2651992439882292
$ws.Range("L5").Value = 0.008197892624240609
$ws.Range("M5").Value = 0.2098692660742133
$ws.Range("N5").Value = 14.60103762125922

$ws.Range("C6").Value = 22.24
$ws.Range("D6").Value = 88.795
$ws.Range("E6").Value = 1.94228338
$ws.Range("F6").Value = 0.06854052000000001
$ws.Range("G6").Value = 0.38291142
$ws.Range("H6").Value = 8.715522399999998
$ws.Range("I6").Value = 3.092565919097597
$ws.Range("J6").Value = 12.34639252609475
$ws.Range("K6").Value = 0.2957574351606878
$ws.Range("L6").Value = 0.01828849026589691
$ws.Range("M6").Value = 0.120800694337972
$ws.Range("N6").Value = 3.54874706657617

$ws.Range("C7").Value = 26.814
$ws.Range("D7").Value = 99.407
$ws.Range("E7").Value = 1.73407214
$ws.Range("F7").Value = 0.05480926
$ws.Range("G7").Value = 0.36421596
$ws.Range("H7").Value = 9.98290697
$ws.Range("I7").Value = 4.789816942884806
$ws.Range("J7").Value = 14.07115305681622
$ws.Range("K7").Value = 0.2453358278238921
$ws.Range("L7").Value = 0.01433519126349422
$ws.Range("M7").Value = 0.1041087796070517
$ws.Range("N7").Value = 4.185612811903998

$ws.Range("C8").Value = 15.228
$ws.Range("D8").Value = 91.149
$ws.Range("E8").Value = 1.94890236
$ws.Range("F8").Value = 0.07774101
$ws.Range("G8").Value = 0.20142543
$ws.Range("H8").Value = 3.25956092
$ws.Range("I8").Value = 3.142911690438326
$ws.Range("J8").Value = 18.84053515070007
$ws.Range("K8").Value = 0.4849272918242346
$ws.Range("L8").Value = 0.02481443004823038
$ws.Range("M8").Value = 0.08846879909905599
$ws.Range("N8").Value = 1.963410230337498

$ws.Range("C9").Value = 18.938
$ws.Range("D9").Value = 98.31
$ws.Range("E9").Value = 1.76433569
$ws.Range("F9").Value = 0.05879493
$ws.Range("G9").Value = 0.18485576
$ws.Range("H9").Value = 3.64024715
$ws.Range("I9").Value = 3.982966233853862
$ws.Range("J9").Value = 15.72339419861771
$ws.Range("K9").Value = 0.2898787217490217
$ws.Range("L9").Value = 0.01605943959779716
$ws.Range("M9").Value = 0.0641082161676848
$ws.Range("N9").Value = 1.950845626872322

$ws.Range("C10").Value = 11.196
$ws.Range("D10").Value = 89.36
$ws.Range("E10").Value = 2.01692432
$ws.Range("F10").Value = 0.07241047
$ws.Range("G10").Value = 0.10395311
$ws.Range("H10").Value = 1.25968689
$ws.Range("I10").Value = 2.678200283167179
$ws.Range("J10").Value = 21.43501480894893
$ws.Range("K10").Value = 0.5478882585675786
$ws.Range("L10").Value = 0.02225214462966929
$ws.Range("M10").Value = 0.04771749559033157
$ws.Range("N10").Value = 0.822380439596942

$ws.Range("C11").Value = 14.61
$ws.Range("D11").Value = 92.39
$ws.Range("E11").Value = 1.8850577
$ws.Range("F11").Value = 0.0589773
$ws.Range("G11").Value = 0.10737908
$ws.Range("H11").Value = 1.64997497
$ws.Range("I11").Value = 3.426607891434566
$ws.Range("J11").Value = 16.04996889409763
$ws.Range("K11").Value = 0.3316769131221047
$ws.Range("L11").Value = 0.01611278534461257
$ws.Range("M11").Value = 0.03924565396714531
$ws.Range("N11").Value = 0.9841689665861194

$ws.Range("C12").Value = 8.502
$ws.Range("D12").Value = 84.763
$ws.Range("E12").Value = 2.17556602
$ws.Range("F12").Value = 0.07067545
$ws.Range("G12").Value = 0.06265933
$ws.Range("H12").Value = 0.59896682
$ws.Range("I12").Value = 2.414917768336601
$ws.Range("J12").Value = 24.08499362503378
$ws.Range("K12").Value = 0.6820477008218383
$ws.Range("L12").Value = 0.02580793968041711
$ws.Range("M12").Value = 0.03628502542263159
$ws.Range("N12").Value = 0.5124487553724474

$ws.Range("C13").Value = 11.681
$ws.Range("D13").Value = 83.903
$ws.Range("E13").Value = 2.09356865
$ws.Range("F13").Value = 0.05140758
$ws.Range("G13").Value = 0.05981359999999999
$ws.Range("H13").Value = 0.7497527900000001
$ws.Range("I13").Value = 3.214276056994883
$ws.Range("J13").Value = 16.6838227316446
$ws.Range("K13").Value = 0.4138599237686337
$ws.Range("L13").Value = 0.01378148061536029
$ws.Range("M13").Value = 0.0232657433677597
$ws.Range("N13").Value = 0.5121972731318662
